# disability_prevalence.xlsx (Zestaponi) - refresh the social-package
# disability table with the "Unified database of targeted social
# assistance program" figures and a revised title / row labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B", "C", "D", "E", "F", "G", "H", "I")

# ---------------------------------------------------------------------
# Insert a new row above the old data row (old row 4, "Number of
# disability persons"). After this:
#   row4 = brand-new empty row (formatting copied from row5 below)
#   row5 = old row4  ("Number of disability persons" + its data)
#   row6 = old row5  (merged source note)
# ---------------------------------------------------------------------
$ws.Rows(4).Insert()

# Seed row 4's look from row 5 (the row that used to be row 4), so the
# new row starts with the same font/fill/wrap as the data row above it.
$ws.Range("A5:I5").Copy()
$ws.Range("A4:I4").PasteSpecial(-4122)

# Use the (untouched) numeric-cell format from B5 (= old B4) as the
# canonical template for every plain data cell in the table (no
# border, no explicit alignment, "#  ##0" number format).
$ws.Range("B5").Copy()
$ws.Range("B4:I4").PasteSpecial(-4122)
$ws.Range("B5").Copy()
$ws.Range("B5:I5").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Row 1 - merged title
# ---------------------------------------------------------------------
$ws.Range("A1:I1").Merge()
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Zestaponi Municipality"
$ws.Range("A1:I1").HorizontalAlignment = -4108
$ws.Range("A1:I1").VerticalAlignment = -4108
$ws.Range("A1:I1").WrapText = $true
$ws.Rows(1).RowHeight = 51

# ---------------------------------------------------------------------
# Row 3 - A3 switches to the Sylfaen font (years header row untouched)
# ---------------------------------------------------------------------
$ws.Range("A3").Font.Name = "Sylfaen"
$ws.Range("A3").Font.Size = 11

# ---------------------------------------------------------------------
# Row 4 - "family with disabilities Persons" (keep the inherited top
# border, drop the bottom one)
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "family with disabilities Persons "
$ws.Range("A4").Borders.Item(9).LineStyle = -4142

$data4 = @(1373, 1300, 1206, 1236, 1211, 1269, 1230, 1217)
for ($i = 0; $i -lt 8; $i++) {
    $ws.Range($cols[$i] + "4").Value = $data4[$i]
}
$ws.Rows(4).RowHeight = 24.75

# ---------------------------------------------------------------------
# Row 5 - "disabilities Persons" (keep the inherited bottom border,
# drop the top one; last cell I5 keeps its own bottom border)
# ---------------------------------------------------------------------
$ws.Range("A5").Value = "disabilities Persons "
$ws.Range("A5").Borders.Item(8).LineStyle = -4142

$data5 = @(1556, 1474, 1359, 1381, 1344, 1402, 1360, 1345)
for ($i = 0; $i -lt 8; $i++) {
    $cell = $ws.Range($cols[$i] + "5")
    $cell.Value = $data5[$i]
    $cell.Borders.Item(8).LineStyle = -4142
}
$ws.Range("I5").Borders.Item(9).LineStyle = 1
$ws.Range("I5").Borders.Item(9).Weight = 2
$ws.Rows(5).RowHeight = 21

# ---------------------------------------------------------------------
# Row 6 - source note (unchanged text, still merged A6:H6); only the
# first cell's top border is removed (continuation cells keep theirs)
# ---------------------------------------------------------------------
$ws.Range("A6").Borders.Item(8).LineStyle = -4142
$ws.Rows(6).RowHeight = 27.75

# ---------------------------------------------------------------------
# Column A width
# ---------------------------------------------------------------------
$ws.Columns(1).ColumnWidth = 20.81640625

$ws.Range("A1").Select()
